$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1757.2
$ws.Range("I19").Value = 2069.7856
$ws.Range("J19").Value = 1027.8334
$ws.Range("K19").Value = 2069.7856
$ws.Range("L19").Value = 1027.8334
$ws.Range("M19").Value = -1894.7856
$ws.Range("N19").Value = -1377.8334

$ws.Range("H34").Value = 1000
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -797
$ws.Range("N34").Value = $null

$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -285
$ws.Range("N36").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8762.969999999999
$ws.Range("I32").Value = 6773.18
$ws.Range("J32").Value = 14981.0625
$ws.Range("K32").Value = 6773.18
$ws.Range("L32").Value = 14981.0625
$ws.Range("M32").Value = -6486.18
$ws.Range("N32").Value = -15555.0625

$ws.Range("H88").Value = 2678.5557
$ws.Range("I88").Value = 3266.6667
$ws.Range("J88").Value = 2384.5
$ws.Range("K88").Value = 3266.6667
$ws.Range("L88").Value = 2384.5
$ws.Range("M88").Value = -2860.6667
$ws.Range("N88").Value = -3196.5

$ws.Range("H91").Value = 2678.5557
$ws.Range("I91").Value = 3266.6667
$ws.Range("J91").Value = 2384.5
$ws.Range("K91").Value = 3266.6667
$ws.Range("L91").Value = 2384.5
$ws.Range("M91").Value = -1862.6667
$ws.Range("N91").Value = -5192.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 54
$ws.Range("I5").Value = 54
$ws.Range("K5").Value = 54
$ws.Range("M5").Value = 59

$ws.Range("H19").Value = 8500
$ws.Range("J19").Value = 8500
$ws.Range("L19").Value = 8500
$ws.Range("N19").Value = -8846

$ws.Range("H54").Value = 1990.75
$ws.Range("I54").Value = 1990.75
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1990.75
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1506.75
$ws.Range("N54").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1178.5135
$ws.Range("I31").Value = 1124.6129
$ws.Range("J31").Value = 1457
$ws.Range("K31").Value = 1124.6129
$ws.Range("L31").Value = 1457
$ws.Range("M31").Value = -829.6129000000001
$ws.Range("N31").Value = -2047

$ws.Range("H34").Value = 1178.5135
$ws.Range("I34").Value = 1124.6129
$ws.Range("J34").Value = 1457
$ws.Range("K34").Value = 1124.6129
$ws.Range("L34").Value = 1457
$ws.Range("M34").Value = -922.6129000000001
$ws.Range("N34").Value = -1861

$ws.Range("H48").Value = 25000
$ws.Range("J48").Value = 25000
$ws.Range("L48").Value = 25000
$ws.Range("N48").Value = -25952

$ws.Range("H132").Value = 3156.7144
$ws.Range("I132").Value = 2846.7144
$ws.Range("J132").Value = 3466.7144
$ws.Range("K132").Value = 8540.143199999999
$ws.Range("L132").Value = 10400.1432
$ws.Range("M132").Value = -6010.143199999999
$ws.Range("N132").Value = -15460.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3025
$ws.Range("I110").Value = 3025
$ws.Range("K110").Value = 9075
$ws.Range("M110").Value = -4985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4622.516
$ws.Range("I122").Value = 4244.5
$ws.Range("J122").Value = 5309.8184
$ws.Range("K122").Value = 12733.5
$ws.Range("L122").Value = 15929.4552
$ws.Range("M122").Value = -10283.5
$ws.Range("N122").Value = -20829.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1965.4736
$ws.Range("I16").Value = 2154.0588
$ws.Range("J16").Value = 362.5
$ws.Range("K16").Value = 2154.0588
$ws.Range("L16").Value = 362.5
$ws.Range("M16").Value = -1984.0588
$ws.Range("N16").Value = -702.5

$ws.Range("H22").Value = 381.0909
$ws.Range("I22").Value = 338
$ws.Range("J22").Value = 417
$ws.Range("K22").Value = 338
$ws.Range("L22").Value = 417
$ws.Range("M22").Value = -43
$ws.Range("N22").Value = -1007

$ws.Range("H27").Value = 381.0909
$ws.Range("I27").Value = 338
$ws.Range("J27").Value = 417
$ws.Range("K27").Value = 338
$ws.Range("L27").Value = 417
$ws.Range("M27").Value = -231
$ws.Range("N27").Value = -631

$ws.Range("H54").Value = 12998.667
$ws.Range("J54").Value = 12998.667
$ws.Range("L54").Value = 12998.667
$ws.Range("N54").Value = -14286.667

$ws.Range("H82").Value = 1670.9375
$ws.Range("I82").Value = 1533.5555
$ws.Range("J82").Value = 1847.5714
$ws.Range("K82").Value = 1533.5555
$ws.Range("L82").Value = 1847.5714
$ws.Range("M82").Value = -1172.5555
$ws.Range("N82").Value = -2569.5714

$ws.Range("H85").Value = 1670.9375
$ws.Range("I85").Value = 1533.5555
$ws.Range("J85").Value = 1847.5714
$ws.Range("K85").Value = 1533.5555
$ws.Range("L85").Value = 1847.5714
$ws.Range("M85").Value = -285.5554999999999
$ws.Range("N85").Value = -4343.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1070
$ws.Range("I13").Value = 1505
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 1505
$ws.Range("L13").Value = 200
$ws.Range("M13").Value = -1365
$ws.Range("N13").Value = -480

$ws.Range("H21").Value = 15153.4
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 21219.143
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 21219.143
$ws.Range("M21").Value = -765
$ws.Range("N21").Value = -21689.143

$ws.Range("H22").Value = 3000
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3586

$ws.Range("H35").Value = 15153.4
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 21219.143
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 21219.143
$ws.Range("M35").Value = -710
$ws.Range("N35").Value = -21799.143

$ws.Range("H47").Value = 14998
$ws.Range("J47").Value = 14998
$ws.Range("L47").Value = 14998
$ws.Range("N47").Value = -16142

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null

$ws.Range("H132").Value = 6633.8335
$ws.Range("I132").Value = 7951
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 23853
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -21323
$ws.Range("N132").Value = -17058.5

